# Updates cryptos list values (Price and Volume(1h) columns) and row reorder for ARBITRUM/FraxShare
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column keeps its original text formatting (values such as "1.00" or
# "43.928.68" must remain literal text, not be re-interpreted as numbers/dates by Excel).
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '43.928.68'
$ws.Range("E2").Value = '  +0.09%  '

# Row 3
$ws.Range("D3").Value = '2.275.26'
$ws.Range("E3").Value = '  +2.62%  '

# Row 4
$ws.Range("E4").Value = '  +0.10%  '

# Row 5
$ws.Range("D5").Value = '269.34'
$ws.Range("E5").Value = '  +3.14%  '

# Row 6
$ws.Range("D6").Value = '94.80'
$ws.Range("E6").Value = '  +9.80%  '

# Row 7
$ws.Range("D7").Value = '0.626'
$ws.Range("E7").Value = '  +1.01%  '

# Row 8
$ws.Range("E8").Value = '  -0.05%  '

# Row 9
$ws.Range("D9").Value = '0.625'
$ws.Range("E9").Value = '  +3.14%  '

# Row 10
$ws.Range("D10").Value = '46.80'
$ws.Range("E10").Value = '  +3.32%  '

# Row 11
$ws.Range("D11").Value = '0.0936'
$ws.Range("E11").Value = '  +1.65%  '

# Row 12
$ws.Range("D12").Value = '8.07'
$ws.Range("E12").Value = '  +7.70%  '

# Row 13
$ws.Range("E13").Value = '  +0.71%  '

# Row 14
$ws.Range("D14").Value = '2.614.12'
$ws.Range("E14").Value = '  +2.50%  '

# Row 15
$ws.Range("D15").Value = '15.41'
$ws.Range("E15").Value = '  +6.50%  '

# Row 16
$ws.Range("D16").Value = '0.830'
$ws.Range("E16").Value = '  +6.54%  '

# Row 17
$ws.Range("D17").Value = '2.275.48'
$ws.Range("E17").Value = '  +2.72%  '

# Row 18
$ws.Range("D18").Value = '43.982.08'
$ws.Range("E18").Value = '  +0.35%  '

# Row 19
$ws.Range("E19").Value = '  +1.77%  '

# Row 20
$ws.Range("D20").Value = '6.18'
$ws.Range("E20").Value = '  +3.94%  '

# Row 21
$ws.Range("D21").Value = '71.07'
$ws.Range("E21").Value = '  +1.66%  '

# Row 22
$ws.Range("D22").Value = '2.30'
$ws.Range("E22").Value = '  -2.01%  '

# Row 23
$ws.Range("D23").Value = '10.03'
$ws.Range("E23").Value = '  +11.79%  '

# Row 24
$ws.Range("D24").Value = '236.53'
$ws.Range("E24").Value = '  +2.67%  '

# Row 25
$ws.Range("D25").Value = '1.00'
$ws.Range("E25").Value = '  -0.01%  '

# Row 26
$ws.Range("D26").Value = '11.37'
$ws.Range("E26").Value = '  +6.77%  '

# Row 27
$ws.Range("D27").Value = '2.49'
$ws.Range("E27").Value = '  +9.89%  '

# Row 28
$ws.Range("D28").Value = '39.49'
$ws.Range("E28").Value = '  -1.61%  '

# Row 29
$ws.Range("E29").Value = '  -5.24%  '

# Row 30
$ws.Range("E30").Value = '  +0.07%  '

# Row 31
$ws.Range("D31").Value = '22.16'
$ws.Range("E31").Value = '  +8.14%  '

# Row 32
$ws.Range("D32").Value = '173.45'
$ws.Range("E32").Value = '  -0.55%  '

# Row 33
$ws.Range("D33").Value = '0.0910'
$ws.Range("E33").Value = '  +4.79%  '

# Row 34
$ws.Range("D34").Value = '5.60'
$ws.Range("E34").Value = '  +4.12%  '

# Row 35
$ws.Range("E35").Value = '  +1.39%  '

# Row 36
$ws.Range("D36").Value = '0.112'
$ws.Range("E36").Value = '  -0.27%  '

# Row 37
$ws.Range("D37").Value = '4.44'
$ws.Range("E37").Value = '  -1.36%  '

# Row 38
$ws.Range("D38").Value = '0.0351'
$ws.Range("E38").Value = '  -1.45%  '

# Row 39
$ws.Range("D39").Value = '3.46'
$ws.Range("E39").Value = '  +18.26%  '

# Row 40
$ws.Range("D40").Value = '0.249'
$ws.Range("E40").Value = '  +24.23%  '

# Row 41
$ws.Range("D41").Value = '2.23'
$ws.Range("E41").Value = '  +6.25%  '

# Row 42
$ws.Range("E42").Value = '  -2.70%  '

# Row 43
$ws.Range("D43").Value = '5.48'
$ws.Range("E43").Value = '  -0.73%  '

# Row 44
$ws.Range("D44").Value = '61.76'
$ws.Range("E44").Value = '  -2.52%  '

# Row 45
$ws.Range("E45").Value = '  +4.62%  '

# Row 46
$ws.Range("B46").Value = 'ARBITRUM'
$ws.Range("C46").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D46").Value = '1.22'
$ws.Range("E46").Value = '  +7.92%  '

# Row 47
$ws.Range("B47").Value = 'FraxShare'
$ws.Range("C47").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D47").Value = '8.53'
$ws.Range("E47").Value = '  +2.33%  '

# Row 48
$ws.Range("D48").Value = '99.65'
$ws.Range("E48").Value = '  -1.09%  '

# Row 49
$ws.Range("D49").Value = '1.19'
$ws.Range("E49").Value = '  +0.59%  '

# Row 50
$ws.Range("D50").Value = '0.430'
$ws.Range("E50").Value = '  -1.86%  '

# Row 51
$ws.Range("D51").Value = '2.491.14'
$ws.Range("E51").Value = '  +2.31%  '

Write-Host "Updated cryptos list values"
